$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 53 and 54 had their match data (columns F:V) swapped - the fixture
#    "St. Pauli v Holstein Kiel" and "Hertha Berlin v Braunschweig" traded
#    places (index/date columns A:E are identical for both rows, so they are
#    left untouched).
# ---------------------------------------------------------------------------
$row53 = @("Hertha Berlin", 3, "Braunschweig", 0, 1.59, "04/09/2023 08:42", 1.78, "17/09/2023 13:20", 4.48, "04/09/2023 08:42", 4.2, "17/09/2023 13:29", 5.2, "04/09/2023 08:42", 4.39, "17/09/2023 13:23", "https://www.betexplorer.com/football/germany/2-bundesliga/hertha-berlin-braunschweig/6PmoIlmf/")

$row54 = @("St. Pauli", 5, "Holstein Kiel", 1, 1.65, "04/09/2023 08:42", 1.93, "17/09/2023 13:29", 4.28, "04/09/2023 08:42", 3.7, "17/09/2023 13:29", 4.96, "04/09/2023 08:42", 4.12, "17/09/2023 13:28", "https://www.betexplorer.com/football/germany/2-bundesliga/st-pauli-holstein-kiel/ny5rH820/")

for ($i = 0; $i -lt $row53.Length; $i++) {
    $ws.Cells.Item(53, 6 + $i).Value = $row53[$i]
}
for ($i = 0; $i -lt $row54.Length; $i++) {
    $ws.Cells.Item(54, 6 + $i).Value = $row54[$i]
}

# ---------------------------------------------------------------------------
# 2) Three new match rows were appended at the bottom of the sheet (rows
#    107-109), extending the used range from A1:V106 to A1:V109. Copy the
#    formatting of the last existing row down onto the new rows first so the
#    index column (bold/bordered) and date column (custom number format)
#    keep their styling, then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A106:V106").Copy()
$ws.Range("A107:V109").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(106, "germany", "2-bundesliga", "2023-2024", 45235.5625, "Hansa Rostock", 0, "Hertha Berlin", 0, 2.66, "29/10/2023 13:42", 2.86, "05/11/2023 13:23", 3.5, "29/10/2023 13:42", 3.58, "05/11/2023 13:26", 2.7, "29/10/2023 13:42", 2.51, "05/11/2023 13:23", "https://www.betexplorer.com/football/germany/2-bundesliga/hansa-rostock-hertha-berlin/0rVncjSQ/"),
    @(107, "germany", "2-bundesliga", "2023-2024", 45235.5625, "Hannover", 2, "Braunschweig", 0, 1.58, "29/10/2023 13:42", 1.55, "05/11/2023 13:04", 4.59, "29/10/2023 13:42", 4.36, "05/11/2023 13:28", 5.16, "29/10/2023 13:42", 6.53, "05/11/2023 13:28", "https://www.betexplorer.com/football/germany/2-bundesliga/hannover-braunschweig/nuRrbWCK/"),
    @(108, "germany", "2-bundesliga", "2023-2024", 45235.5625, "Karlsruher SC", 0, "Paderborn", 3, 2.01, "29/10/2023 13:42", 2.23, "05/11/2023 13:04", 4.05, "29/10/2023 13:42", 3.68, "05/11/2023 12:58", 3.37, "29/10/2023 13:42", 3.24, "05/11/2023 12:56", "https://www.betexplorer.com/football/germany/2-bundesliga/karlsruher-paderborn/QPxz0hr8/")
)

$r = 107
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
    $ws.Cells.Item($r, 16).Value = $row[15]
    $ws.Cells.Item($r, 17).Value = $row[16]
    $ws.Cells.Item($r, 18).Value = $row[17]
    $ws.Cells.Item($r, 19).Value = $row[18]
    $ws.Cells.Item($r, 20).Value = $row[19]
    $ws.Cells.Item($r, 21).Value = $row[20]
    $ws.Cells.Item($r, 22).Value = $row[21]
    $r++
}
